# Updates cryptos list values per diff (Price and Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.549.46"
$ws.Range("E2").Value = "  +1.40%  "
$ws.Range("D3").Value = "3.950.27"
$ws.Range("E3").Value = "  +4.62%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'479.42"
$ws.Range("E5").Value = "  +9.13%  "
$ws.Range("D6").Value = "'148.87"
$ws.Range("E6").Value = "  +4.20%  "
$ws.Range("D7").Value = "'0.628"
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("D9").Value = "'0.734"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("E10").Value = "  +11.31%  "
$ws.Range("D11").Value = "'0.0000353"
$ws.Range("E11").Value = "  +14.24%  "
$ws.Range("D12").Value = "'43.47"
$ws.Range("D13").Value = "4.578.12"
$ws.Range("E13").Value = "  +4.25%  "
$ws.Range("D14").Value = "'10.47"
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("D15").Value = "'15.03"
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("D16").Value = "3.946.51"
$ws.Range("E16").Value = "  +4.08%  "
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "'20.11"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("E19").Value = "  +2.45%  "
$ws.Range("D20").Value = "67.724.56"
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("D21").Value = "'435.09"
$ws.Range("E21").Value = "  +6.03%  "
$ws.Range("E22").Value = "  +4.13%  "
$ws.Range("D23").Value = "'14.55"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "'87.62"
$ws.Range("E24").Value = "  +2.86%  "
$ws.Range("E25").Value = "  +8.17%  "
$ws.Range("D26").Value = "'38.73"
$ws.Range("E26").Value = "  +5.34%  "
$ws.Range("D27").Value = "'10.18"
$ws.Range("E27").Value = "  +5.01%  "
$ws.Range("D28").Value = "'9.80"
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("D29").Value = "'718.74"
$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("D31").Value = "'13.48"
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("E32").Value = "  +4.78%  "
$ws.Range("D33").Value = "'42.30"
$ws.Range("E33").Value = "  -0.99%  "
$ws.Range("D34").Value = "0.0₃0848"
$ws.Range("E34").Value = "  +28.45%  "
$ws.Range("D35").Value = "'58.34"
$ws.Range("E35").Value = "  +3.86%  "
$ws.Range("E36").Value = "  -3.12%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("E38").Value = "  -2.93%  "
$ws.Range("D39").Value = "'0.0477"
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("D40").Value = "'3.08"
$ws.Range("E40").Value = "  +6.48%  "
$ws.Range("E41").Value = "  +3.35%  "
$ws.Range("E42").Value = "  +2.90%  "
$ws.Range("E43").Value = "  +7.35%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "'2.83"
$ws.Range("E45").Value = "  +6.26%  "
$ws.Range("E46").Value = "  +5.68%  "
$ws.Range("D47").Value = "'2.52"
$ws.Range("E47").Value = "  -6.85%  "
$ws.Range("D48").Value = "'3.25"
$ws.Range("E48").Value = "  -2.39%  "
$ws.Range("D49").Value = "'149.83"
$ws.Range("E49").Value = "  +4.77%  "
$ws.Range("E50").Value = "  +3.02%  "
$ws.Range("D51").Value = "'25.63"
$ws.Range("E51").Value = "  +3.82%  "
